# Preference Flow.xlsx - update preference-flow formulas.
#
# 1) J6:J12 "B-minus-next-row" deltas are re-entered as one relative fill
#    across the range (mirrors Excel turning the individually authored
#    formulas B6-B7, B7-B8, ... into one filled/shared formula block).
# 2) The RMSE / Count helper formulas in O15, O16, O19 and O20 are widened
#    to also take in the K5:L12 / K6:L12 preference-flow delta ranges,
#    so the UAP (and other newly-tracked) columns feed into the stats.
# 3) Leave the cursor on J3, matching where the author ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-fill J6:J12 with the same relative formula (B_row - B_row+1) ---
$ws.Range("J6:J12").FormulaR1C1 = "=RC[-8]-R[1]C[-8]"

# --- 2) Widen the RMSE/Count aggregate formulas to include K:L deltas ---
$ws.Range("O15").Formula = "=SQRT(SUMSQ(B22:K25,K5:L12)/COUNT(B22:K25,K5:L12))"
$ws.Range("O16").Formula = "=COUNT(B22:K25,K5:L12)"
$ws.Range("O19").Formula = "=SQRT(SUMSQ(B23:K25,K6:L12)/COUNT(B23:K25,K6:L12))"
$ws.Range("O20").Formula = "=COUNT(B23:K25,K6:L12)"

# --- 3) Update the saved selection/active cell to J3 ---
$ws.Range("J3").Select()
